# SE306-FinalChecklist.docx edit script
# Applies the changes described by the commit "Updated checklist and the boss cutscene"
#
# Strategy: most hunks in the diff simply remove w:proofErr (spell/grammar
# check) markers and merge the runs they used to separate. Doing a
# Find/Replace over a range that spans such a proofErr (i.e. the match
# starts with real text BEFORE the proofErr and continues past it)
# naturally causes the engine to rebuild that span as a single run with
# no proofErr markers left over - exactly matching the target XML. Where
# the visible text itself changes, the Find/Replace performs that text
# substitution at the same time.
#
# A handful of proofErr-wrapped words are the very first thing in their
# paragraph/cell (nothing precedes them), so a plain Find/Replace leaves
# the proofErr behind (nothing "before" it is in the matched range). For
# those we temporarily insert one throw-away character before the
# paragraph, include it in the Find string, and have the Replace drop it
# again - this forces the whole span (including the leading proofErr) to
# be rebuilt as a single clean run.

$d = $word.ActiveDocument

function Replace-Text($find, $replace, $matchCase = $true) {
    $rng = $d.Content
    $rng.Find.Execute($find, $matchCase, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

function Fix-LeadingProofErr($paraRange, $oldFullText, $newFullText, $matchCase = $true) {
    $insPt = $paraRange.Duplicate
    $insPt.Collapse(1)  # wdCollapseStart
    $insPt.InsertBefore("\")
    Replace-Text ("\" + $oldFullText) $newFullText $matchCase
}

# 1. "...your initial project plan to [gramStart]completed[gramEnd] by the final deadline."
#    -> merge runs, remove gramStart/gramEnd proofErr, no text change.
Replace-Text "your initial project plan to completed by the final deadline." "your initial project plan to completed by the final deadline."

# 2. Signature area "___Prisma___" -> remove spellStart/spellEnd around Prisma (no text change)
Replace-Text "_________________Prisma________________________" "_________________Prisma________________________"

# 3. Four "Prisma features ..." table cells -> merge away spellStart/spellEnd (no text change).
#    Prisma is the first word in each of these cells, so use the leading-proofErr fix.
$t1 = $d.Tables.Item(1)
Fix-LeadingProofErr $t1.Cell(3, 5).Range.Paragraphs.Item(1).Range `
    "Prisma features a unique main character which functions as expected." `
    "Prisma features a unique main character which functions as expected."
Fix-LeadingProofErr $t1.Cell(5, 5).Range.Paragraphs.Item(1).Range `
    "Prisma features 6 levels, each with a unique design " `
    "Prisma features 6 levels, each with a unique design "
Fix-LeadingProofErr $t1.Cell(9, 5).Range.Paragraphs.Item(1).Range `
    "Prisma features a welcome screen which allows the player to navigate to the stage select, " `
    "Prisma features a welcome screen which allows the player to navigate to the stage select, "
Fix-LeadingProofErr $t1.Cell(11, 5).Range.Paragraphs.Item(1).Range `
    "Prisma features 6 levels of varying complexity. " `
    "Prisma features 6 levels of varying complexity. "

# 4. "...ending cutscene. " -> merge away spellStart/spellEnd (no text change)
Replace-Text "with the end signified by an ending cutscene. " "with the end signified by an ending cutscene. "

# 5. "...documented through cutscenes which..." -> merge away spellStart/spellEnd (no text change)
Replace-Text "The game’s storyline is documented through cutscenes which are unlocked as the player progresses." "The game’s storyline is documented through cutscenes which are unlocked as the player progresses."

# 6. "The game has been playtested thoroughly within the team," -> merge away spellStart/spellEnd
Replace-Text "The game has been playtested thoroughly within the team," "The game has been playtested thoroughly within the team,"

# 7. "Team UGSoft playtested our game..." -> merge away spellStart/spellEnd x2 and drop _GoBack bookmark
Replace-Text "Team UGSoft playtested our game and suggested how to improve it" "Team UGSoft playtested our game and suggested how to improve it"

# 8. "The game’s UI changes color according to their current progress in the game." -> merge away spellStart/spellEnd
Replace-Text "The game’s UI changes color according to their current progress in the game." "The game’s UI changes color according to their current progress in the game."

# 9. "+ Leaderboard" (two occurrences) -> merge away spellStart/spellEnd around Leaderboard
Replace-Text "Local Multiplayer + Leaderboard.  10%" "Local Multiplayer + Leaderboard.  10%"
Replace-Text "Online Multiplayer + Leaderboard. 10%" "Online Multiplayer + Leaderboard. 10%"

# 10. Add "X" markers to the 6 empty cells in the Local/Online Multiplayer rows
$t1.Cell(20, 2).Range.Text = "X"
$t1.Cell(20, 3).Range.Text = "X"
$t1.Cell(20, 4).Range.Text = "X"
$t1.Cell(21, 2).Range.Text = "X"
$t1.Cell(21, 3).Range.Text = "X"
$t1.Cell(21, 4).Range.Text = "X"

# 11. "Monetization options have been discussed within the Github’s wiki page" -> merge away spellStart/spellEnd
Replace-Text "Monetization options have been discussed within the Github’s wiki page" "Monetization options have been discussed within the Github’s wiki page"

# 12. "Color-blind mode. Claiming 7%" -> "Color-blind mode. Claiming 4%".
#     "Color" is the first word of this cell, so use the leading-proofErr fix.
Fix-LeadingProofErr $t1.Cell(29, 1).Range.Paragraphs.Item(1).Range `
    "Color-blind mode. Claiming 7%" `
    "Color-blind mode. Claiming 4%"

# 13. "Due to Prisma’s heavy reliance on color ... settings menu." -> merge away all the spellStart/spellEnd markers
Replace-Text "Due to Prisma’s heavy reliance on color throughout the game, color blind modes have been implemented, and are able to be set in the settings menu." "Due to Prisma’s heavy reliance on color throughout the game, color blind modes have been implemented, and are able to be set in the settings menu."

# 14. "...progress to the cloud – Describe here. Claiming 3%" -> "...Claiming 6%"
Replace-Text "progress to the cloud – Describe here. Claiming 3%" "progress to the cloud – Describe here. Claiming 6%"

# 15. "Levels are completable" -> merge away spellStart/spellEnd
Replace-Text "Levels are completable" "Levels are completable"

# 16. "Each level has been playtested to completion." -> merge away spellStart/spellEnd
Replace-Text "Each level has been playtested to completion." "Each level has been playtested to completion."

# 17. "The game has been playtested with those identifying within the user-group, and met with positive reactions." -> merge away spellStart/spellEnd
Replace-Text "The game has been playtested with those identifying within the user-group, and met with positive reactions." "The game has been playtested with those identifying within the user-group, and met with positive reactions."

# 18. Append sentence to "The work was distributed fairly and according to people’s strengths,"
Replace-Text "The work was distributed fairly and according to people’s strengths," "The work was distributed fairly and according to people’s strengths, and / or willingness to learn about a certain feature."

# 19. "Documentation regrading how the project..." -> "Documentation regarding how the project..." (also merges away spellStart/spellEnd)
Replace-Text "Documentation regrading how the project was managed and matters relating to team work were documented in the wiki." "Documentation regarding how the project was managed and matters relating to team work were documented in the wiki."

# 20. Move the page-break split point in the meeting-minutes sentence (net text identical)
Replace-Text "Meeting minutes were created at the end of each meeting to ensure members not present were updated. A few were rather uninformative, due to all members being present and a cohesive understanding of what was to be done." "Meeting minutes were created at the end of each meeting to ensure members not present were updated. A few were rather uninformative, due to all members being present and a cohesive understanding of what was to be done."

# 21. "A large amount of the art was hand made using photoshop, ..." -> merge away spellStart/spellEnd
Replace-Text "A large amount of the art was hand made using photoshop, with the rest of the assets either bought or under creative commons licenses. The original assets can be found on the wiki page." "A large amount of the art was hand made using photoshop, with the rest of the assets either bought or under creative commons licenses. The original assets can be found on the wiki page."

# 22. "Design Decisions (SoftEng and Game Design)" -> merge away spellStart/spellEnd
Replace-Text "Design Decisions (SoftEng and Game Design)" "Design Decisions (SoftEng and Game Design)"

# 23. Append sentence to "Information regarding the "
Replace-Text "Information regarding the " "Information regarding the teams thoughts on how we worked together and the RUP process are available on the Github Wiki page"

# 24. "70%" -> "60%" (Extent of Development/Scripting in Unity)
$t2 = $d.Tables.Item(2)
$t2.Cell(28, 2).Range.Text = "60%"

# 25. Rewrite the "Although we used Google..." paragraph
Replace-Text "Although we used Google a copious amount of times to find solutions, we often adjusted the code greatly to fit out needs, forming logic which, although based off " "Although we used Google a copious amount of times to find solutions, we often adjusted the code greatly to fit our needs, forming logic which, although based off solutions found online, was transformed to fit our needs."

# 26. Fill in the final empty cell with resource information
$t2.Cell(29, 3).Range.Text = "The resources generated by the team were the orbs, backgrounds, icons, the title screen assets. More in depth information about our original assets can be found on the github wiki."
